$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Component Placement List (CPL) data - v4
# Columns: Row, Designator, Mid X, Mid Y, Layer, Rotation
$data = @()
$data += ,@(2, "J3", 131.35, -161.14, "top", 180)
$data += ,@(3, "C5", 110.1, -161.2, "top", 90)
$data += ,@(4, "R9", 107.85, -162.16, "top", 0)
$data += ,@(5, "D3", 148.75, -160.1, "top", 180)
$data += ,@(6, "C4", 105.85, -155.175, "top", 0)
$data += ,@(7, "J3_F2_BAT1", 112.2, -170.2, "top", -90)
$data += ,@(8, "C2", 114.6, -161.2, "top", 90)
$data += ,@(9, "J4_M0_MOT0", 101.067857, -168.832144, "top", -135)
$data += ,@(10, "J1", 135.1, -161.14, "top", 180)
$data += ,@(11, "R7", 116.8375, -160.8, "top", 0)
$data += ,@(12, "C3", 111.6, -161.2, "top", 90)
$data += ,@(13, "J4_M3_MOT3", 127.8, -186.75, "top", 180)
$data += ,@(14, "D7", 116.8375, -159.2, "top", 0)
$data += ,@(15, "R8", 104.85, -162.2, "top", 0)
$data += ,@(16, "D1", 148.75, -157.1, "top", 180)
$data += ,@(17, "D5", 148.3, -190.55, "top", 0)
$data += ,@(18, "J2", 76.825, -161.2, "top", 180)
$data += ,@(19, "J4", 127.6, -161.14, "top", 180)
$data += ,@(20, "MAN1", 96.59999999999999, -154.45, "top", 180)
$data += ,@(21, "J9", 100.35, -160.2, "top", 180)
$data += ,@(22, "R3", 76.09999999999999, -188.95, "top", 180)
$data += ,@(23, "R11", 107.85, -160.66, "top", 180)
$data += ,@(24, "R1", 148.35, -146.65, "top", 0)
$data += ,@(25, "D4", 148.35, -145.05, "top", 0)
$data += ,@(26, "GND1", 143.75, -157.4, "top", -90)
$data += ,@(27, "R15", 92.34999999999999, -154.7, "top", -90)
$data += ,@(28, "IMU1", 111.35, -157.7, "top", 180)
$data += ,@(29, "R14", 90.84999999999999, -154.7, "top", -90)
$data += ,@(30, "5V1", 143.75, -161.15, "top", -90)
$data += ,@(31, "R2", 148.3, -188.95, "top", 0)
$data += ,@(32, "RPI2", 88.395, -146.505, "top", 90)
$data += ,@(33, "J4_M2_MOT2", 123.182144, -168.782144, "top", -45)
$data += ,@(34, "R10", 104.85, -160.7, "top", 180)
$data += ,@(35, "D6", 76.09999999999999, -190.55, "top", 180)
$data += ,@(36, "C1", 113.1, -161.2, "top", 90)
$data += ,@(37, "J7", 123.85, -161.14, "top", 180)
$data += ,@(38, "J37", 143.76, -153.65, "top", -90)
$data += ,@(39, "D2", 148.75, -158.6, "top", 180)
$data += ,@(40, "R4", 149.4, -153.05, "top", -90)
$data += ,@(41, "R5", 147.9, -153.05, "top", -90)
$data += ,@(42, "J4_M1_MOT1", 96.45, -186.85, "top", 0)
$data += ,@(43, "R6", 146.4, -153.05, "top", -90)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "Done applying CPL v4 update"